$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row text updates (summer-reporting column header relabeling)
$ws.Range("E1").Value = "Legal Given Name"
$ws.Range("F1").Value = "Birthdate"
$ws.Range("G1").Value = "Ministry Course Code and Level"
$ws.Range("I1").Value = "Final Percent"
$ws.Range("K1").Value = "Credits"

# Selection moved to the header row
$ws.Range("A1:K1").Select()
